$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Negate the existing "salaire_total" values for rows 2-4
$ws.Range("D2").Value = -1500
$ws.Range("D3").Value = -704
$ws.Range("D4").Value = -400

# Add new row 5: Jeanno
$ws.Range("A5").Value = "Jeanno"
$ws.Range("B5").Value = 10000
$ws.Range("C5").Value = 0
$ws.Range("D5").Value = -10000

# Add new row 6: Thai Nhien
$ws.Range("A6").Value = "Thai Nhien"
$ws.Range("B6").Value = 20000
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = -20000
